$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (Summary / Labels / Tests Stories) - shifts data rows up
$ws.Rows.Item(1).Delete()

# Add a new trailing value in column A of the now-last row
$ws.Range("A4").Value = "a"

# Move the active selection to A5 (the first empty row below the data)
$ws.Range("A5").Select() | Out-Null
